$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-06-03 Tuesday" "2025-06-04 Wednesday"

Replace-Text "217×6=" "236×9="
Replace-Text "619×5=" "148×8="
Replace-Text "111×8=" "501×6="
Replace-Text "404×2=" "994×6="
Replace-Text "363×6=" "504×7="
Replace-Text "547×9=" "853×6="
Replace-Text "310×8=" "713×7="
Replace-Text "315×5=" "890×6="
Replace-Text "955×7=" "170×3="
Replace-Text "466×2=" "591×3="
Replace-Text "923×2=" "127×4="
Replace-Text "739×3=" "718×2="
Replace-Text "818×6=" "563×3="
Replace-Text "382×9=" "494×3="
Replace-Text "700×3=" "933×9="
Replace-Text "485×3=" "962×9="
Replace-Text "556×3=" "713×7="
Replace-Text "343×2=" "435×5="
Replace-Text "619×2=" "308×2="
Replace-Text "199×6=" "744×2="
Replace-Text "677×3=" "381×9="
Replace-Text "187×6=" "874×2="
Replace-Text "181×8=" "240×2="
Replace-Text "782×6=" "778×8="
Replace-Text "670×6=" "584×3="
